# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet.
# The previous last row (46) loses its "last row" date format and gets the
# standard datetime format instead; the new last row (47) takes over the
# "last row" date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 is no longer the last row in the table - restore its date cell to
# the standard timestamp number format used by every other data row.
$ws.Range("A46").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 47: today's tallies.
$ws.Range("A47").NumberFormat = "YYYY-MM-DD"
$ws.Range("A47").Value = 45632
$ws.Range("B47").Value = 121
$ws.Range("C47").Value = 105
$ws.Range("D47").Value = 113
